# Update the cryptocurrency price (column D) and 1h-volume-change (column
# E) figures for this run. Column D cells are stored as text in the sheet
# (values like "65.170.19" or "0.0000191" aren't valid numbers), so D
# values are written with a leading apostrophe - the standard Excel
# text-qualifier - to force text entry instead of numeric auto-conversion.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''65.170.19'
$ws.Range('E2').Value = '  +2.46%  '
$ws.Range('D3').Value = '''2.641.76'
$ws.Range('E3').Value = '  +1.13%  '
$ws.Range('D5').Value = '''601.60'
$ws.Range('E5').Value = '  +1.87%  '
$ws.Range('D6').Value = '''156.59'
$ws.Range('E6').Value = '  +4.45%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('E8').Value = '  +0.76%  '
$ws.Range('E9').Value = '  +10.31%  '
$ws.Range('E10').Value = '  +6.18%  '
$ws.Range('D11').Value = '''5.80'
$ws.Range('E11').Value = '  +0.50%  '
$ws.Range('D12').Value = '''0.154'
$ws.Range('E12').Value = '  +2.07%  '
$ws.Range('D13').Value = '''29.37'
$ws.Range('E13').Value = '  +6.06%  '
$ws.Range('D14').Value = '''0.0000191'
$ws.Range('E14').Value = '  +21.59%  '
$ws.Range('D15').Value = '''3.117.07'
$ws.Range('E15').Value = '  +1.39%  '
$ws.Range('D16').Value = '''64.992.75'
$ws.Range('E16').Value = '  +2.48%  '
$ws.Range('D17').Value = '''2.649.66'
$ws.Range('E17').Value = '  +1.56%  '
$ws.Range('D18').Value = '''12.65'
$ws.Range('E18').Value = '  +4.35%  '
$ws.Range('E19').Value = '  +3.45%  '
$ws.Range('D20').Value = '''359.99'
$ws.Range('E20').Value = '  +4.11%  '
$ws.Range('E21').Value = '  +7.20%  '
$ws.Range('E22').Value = '  -0.02%  '
$ws.Range('D23').Value = '''69.26'
$ws.Range('E23').Value = '  +3.45%  '
$ws.Range('E24').Value = '  +0.64%  '
$ws.Range('D25').Value = '''9.43'
$ws.Range('E25').Value = '  +2.18%  '
$ws.Range('D26').Value = '''1.65'
$ws.Range('E26').Value = '  -0.50%  '
$ws.Range('E27').Value = '  -1.72%  '
$ws.Range('E28').Value = '  +2.84%  '
$ws.Range('D29').Value = '''0.0₃0974'
$ws.Range('E29').Value = '  +12.17%  '
$ws.Range('D30').Value = '''549.63'
$ws.Range('E30').Value = '  +0.32%  '
$ws.Range('E31').Value = '  +8.84%  '
$ws.Range('D32').Value = '''0.996'
$ws.Range('E32').Value = '  -0.40%  '
$ws.Range('E33').Value = '  +1.53%  '
$ws.Range('D34').Value = '''5.64'
$ws.Range('E34').Value = '  +5.51%  '
$ws.Range('D35').Value = '''6.37'
$ws.Range('E35').Value = '  +4.06%  '
$ws.Range('E36').Value = '  +4.71%  '
$ws.Range('D37').Value = '''20.50'
$ws.Range('E37').Value = '  +5.09%  '
$ws.Range('D38').Value = '''2.01'
$ws.Range('E38').Value = '  +2.65%  '
$ws.Range('D39').Value = '''162.43'
$ws.Range('E39').Value = '  -2.25%  '
$ws.Range('E40').Value = '  -0.01%  '
$ws.Range('E41').Value = '  +0.05%  '
$ws.Range('D42').Value = '''42.83'
$ws.Range('E42').Value = '  +7.89%  '
$ws.Range('D43').Value = '''166.87'
$ws.Range('E43').Value = '  +0.99%  '
$ws.Range('E44').Value = '  +4.14%  '
$ws.Range('E45').Value = '  +7.26%  '
$ws.Range('E46').Value = '  +7.41%  '
$ws.Range('D47').Value = '''23.26'
$ws.Range('E47').Value = '  +0.76%  '
$ws.Range('D48').Value = '''0.653'
$ws.Range('E48').Value = '  +3.28%  '
$ws.Range('E49').Value = '  +4.95%  '
$ws.Range('E51').Value = '  +2.34%  '
